$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "My order 1324 is delayed, when will it ship?"
$ws.Range("C2").Value = "{'orders': [{'Product ID': 38035, 'ProductName': 'liebherr cufr3311 fire red smartfrost fridge freezer', 'Category': 'Fridge Freezers', 'Category ID': 2622, 'OrderID': 1324, 'CustomerID': 632, 'OrderStatus': 'Delivered', 'ReturnEligible': False, 'ShippingDate': '2024-10-06 12:24:37.765490'}]}"

$ws.Range("B3").Value = "I haven't received my order 1653. Can you tell me when it will arrive?"
$ws.Range("C3").Value = "{'orders': [{'Product ID': 13069, 'ProductName': 'jvc lt 32c460 32 led tv', 'Category': 'TVs', 'Category ID': 2614, 'OrderID': 1653, 'CustomerID': 668, 'OrderStatus': 'Delivered', 'ReturnEligible': True, 'ShippingDate': '2024-10-20 12:24:37.765490'}]}"

$ws.Range("B4").Value = "My order 2713 has been shipped. When will it arrive?"
$ws.Range("C4").Value = "{'orders': [{'Product ID': 12517, 'ProductName': 'samsung ue55mu6500u 6 series 55 led tv', 'Category': 'TVs', 'Category ID': 2614, 'OrderID': 2713, 'CustomerID': 2386, 'OrderStatus': 'Cancelled', 'ReturnEligible': False, 'ShippingDate': '2024-10-11 12:24:37.765490'}]}"

$ws.Range("B5").Value = "The estimated delivery date for order 1306 has passed. Can I get an update?"
$ws.Range("C5").Value = "{'orders': [{'Product ID': 10103, 'ProductName': 'samsung ue75nu8000txzg fernseher premium uhd', 'Category': 'TVs', 'Category ID': 2614, 'OrderID': 1306, 'CustomerID': 571, 'OrderStatus': 'Shipped', 'ReturnEligible': False, 'ShippingDate': '2024-10-22 12:24:37.765490'}]}"

$ws.Range("B6").Value = "I want to cancel order 3368 because it is delayed."
$ws.Range("C6").Value = "{'orders': [{'Product ID': 36534, 'ProductName': 'indesit ui8f1cw freezer', 'Category': 'Freezers', 'Category ID': 2621, 'OrderID': 3368, 'CustomerID': 2073, 'OrderStatus': 'Delivered', 'ReturnEligible': False, 'ShippingDate': '2024-10-28 12:24:37.765490'}]}"

$ws.Range("B7").Value = "I haven't received my order with Order ID 3065 yet. Is it delayed?"
$ws.Range("C7").Value = "{'orders': [{'Product ID': 22065, 'ProductName': 'sony dsc rx0 1.0 type sensor ultra compact camera with waterproof and shockproof design black', 'Category': 'Digital Cameras', 'Category ID': 2617, 'OrderID': 3421, 'CustomerID': 1607, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-27 12:24:37.765490'}]}"

$ws.Range("B8").Value = "The estimated delivery date for order 1985 has passed. When can I expect it?"
$ws.Range("C8").Value = "{'orders': [{'Product ID': 11032, 'ProductName': 'samsung ue65nu7500/ue65nu7500', 'Category': 'TVs', 'Category ID': 2614, 'OrderID': 2640, 'CustomerID': 665, 'OrderStatus': 'Cancelled', 'ReturnEligible': False, 'ShippingDate': '2024-10-26 12:24:37.765490'}]}"

$ws.Range("B9").Value = "My order 1518 is delayed. What's the reason for the delay and when will it be shipped?"
$ws.Range("C9").Value = "{'orders': [{'Product ID': 39112, 'ProductName': 'miele kfn 28132 d freestanding fridge freezer with frost free clean steel', 'Category': 'Fridge Freezers', 'Category ID': 2622, 'OrderID': 2515, 'CustomerID': 1606, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-02 12:24:37.765490'}]}"

$ws.Range("B10").Value = "Order 1874 hasn't arrived yet. Can you give me an update on its shipping status and the new estimated delivery date?"
$ws.Range("C10").Value = "{'orders': [{'Product ID': 46368, 'ProductName': 'smeg fab5rgo refrigerator', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 1631, 'CustomerID': 1917, 'OrderStatus': 'Delivered', 'ReturnEligible': True, 'ShippingDate': '2024-10-04 12:24:37.765490'}]}"

$ws.Range("B11").Value = "Order 3249 is late. I'm very disappointed. Can you tell me when it will arrive?"
$ws.Range("C11").Value = "{'orders': [{'Product ID': 44943, 'ProductName': 'zanussi free standing larder fridge in white', 'Category': 'Fridges', 'Category ID': 2623, 'OrderID': 1459, 'CustomerID': 616, 'OrderStatus': 'Pending', 'ReturnEligible': False, 'ShippingDate': '2024-10-15 12:24:37.765490'}]}"
